$d = $word.ActiveDocument
$author = 'Siemes, Jort (193186)'
$initials = 'JS'

# ------------------------------------------------------------------
# Comment 0 - wraps "Volgens AMC Entertainment ... concertfilms. "
# ------------------------------------------------------------------
$r0a = $d.Content
$r0a.Find.Execute('Volgens AMC Entertainment kwamen de ticketsales uit') | Out-Null
$r0b = $d.Content
$r0b.Find.Execute('Eras Tour film een van de meest bekeken concertfilms. ') | Out-Null
$r0 = $d.Range($r0a.Start, $r0b.End)
$c0 = $d.Comments.Add($r0, "Leuk om hier te vermelden dat ze deze film zelf heeft uitgebracht zonder enige tussenpartijen. Geen studio gebruiken om je film te 'publishen' is namelijk ook een best grote achievement")
$c0.Author = $author
$c0.Initial = $initials

# ------------------------------------------------------------------
# Comment 1 - wraps the whole "In Rio de Janeiro ... concert.  " paragraph
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(32)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$c1 = $d.Comments.Add($r1, 'Ik begreep zelf dat ook het stadion zelf waterflessen meebrengen had verboden? Misschien goed om ook naar te kijken, al wel interessante ontwikkeling rondom "shaduwzijde"')
$c1.Author = $author
$c1.Initial = $initials

# ------------------------------------------------------------------
# Text edit - "privéjets" -> "privéjet" + highlighted "s"
# ------------------------------------------------------------------
$rp = $d.Content
$rp.Find.Execute('privéjets') | Out-Null
$rs = $d.Range($rp.End - 1, $rp.End)
$rs.Find.ClearFormatting()
$rs.Find.Replacement.ClearFormatting()
$rs.Find.Replacement.Highlight = $true
$rs.Find.Execute('s', $false, $false, $false, $false, $false, $true, 1, $false, 's', 2) | Out-Null

# ------------------------------------------------------------------
# Comment 2 - wraps "CO2 "
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute('CO2 ') | Out-Null
$c2 = $d.Comments.Add($r2, 'Ik zou zelf CO2 weglaten, stikstofoxides, roet en andere broeikasgassen worden namelijk door een vliegtuig ook uitgestoten.')
$c2.Author = $author
$c2.Initial = $initials

# ------------------------------------------------------------------
# Comment 3 - wraps "Dit samen zorgt ervoor dat Taylor Swifr niet zo goed is voor het milieu. "
# ------------------------------------------------------------------
$r3a = $d.Content
$r3a.Find.Execute('Dit samen zorgt ervoor dat Taylor') | Out-Null
$r3b = $d.Content
$r3b.Find.Execute('niet zo goed is voor het milieu. ') | Out-Null
$r3 = $d.Range($r3a.Start, $r3b.End)
$c3 = $d.Comments.Add($r3, 'Heel sterk stuk! Alleen vind ik de laatste zin nog wat kracht missen. Ik zou zelf iets meer de koppeling leggen tussen wat deze gigantische populariteit met zich mee sleept. En dat de punten in deze alinea de kosten zijn van de fans & samenleving zijn rondom de tours van deze mega-sterren.')
$c3.Author = $author
$c3.Initial = $initials

# ------------------------------------------------------------------
# Comment 4 - wraps the whole "'Wat is het Taylor Swift Effect ..." paragraph
# ------------------------------------------------------------------
$p4 = $d.Paragraphs(36)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$c4 = $d.Comments.Add($r4, 'Nog wel erg samenvattend, misschien leuk hier de situatie te schetsen voor wanneer Taylor langskomt in Nederland. Sinds je dit in het begin al beetje introduceert ☺️')
$c4.Author = $author
$c4.Initial = $initials

Write-Host "Done. Comments count:" $d.Comments.Count
